$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.167.08'
$ws.Range('E2').Value = '  +5.26%  '
$ws.Range('D3').Value = '2.757.89'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.74'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.21'
$ws.Range('E6').Value = '  +6.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.610'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').Value = '2.758.01'
$ws.Range('E9').Value = '  +2.99%  '
$ws.Range('E10').Value = '  +1.70%  '
$ws.Range('E11').Value = '  +4.54%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.161'
$ws.Range('E12').Value = '  +4.53%  '
$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.388'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = '3.237.67'
$ws.Range('E14').Value = '  +2.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.27'
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('D16').Value = '64.031.56'
$ws.Range('E16').Value = '  +5.04%  '
$ws.Range('E17').Value = '  +5.79%  '
$ws.Range('D18').Value = '2.750.01'
$ws.Range('E18').Value = '  +3.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.93'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('E20').Value = '  +2.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '360.74'
$ws.Range('E21').Value = '  +2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.96'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.531'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.36'
$ws.Range('E25').Value = '  +3.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.170'
$ws.Range('E26').Value = '  +4.81%  '
$ws.Range('E27').Value = '  +4.10%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '0.0₃0914'
$ws.Range('E29').Value = '  +11.62%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.09'
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '171.57'
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('E33').Value = '  +13.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.43'
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.78'
$ws.Range('E36').Value = '  +7.03%  '
$ws.Range('E37').Value = '  +8.14%  '
$ws.Range('E38').Value = '  +9.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.02'
$ws.Range('E39').Value = '  +16.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '345.88'
$ws.Range('E40').Value = '  +5.02%  '
$ws.Range('E41').Value = '  +5.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.44'
$ws.Range('E42').Value = '  +2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.60'
$ws.Range('E43').Value = '  +6.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.74'
$ws.Range('E44').Value = '  +5.68%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0591'
$ws.Range('E45').Value = '  +5.19%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.69'
$ws.Range('E46').Value = '  +5.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.648'
$ws.Range('E47').Value = '  +5.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '139.00'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.01%  '
